$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ProjFunctions": fill in Plan/Actual period numbers for several
# activity rows (editCourse / removeCourse and related class-management
# rows), widen column A by one unit, and mark function 3.10 with a
# 2-decimal display so it reads "3.10" instead of "3.1".
# ---------------------------------------------------------------------
$projFunctions = $wb.Worksheets.Item("ProjFunctions")

function Set-PlanActual($ws, $row, $planStart, $planDur, $actualStart, $actualDur, $percent) {
    $ws.Cells.Item($row, 5).Value = $planStart
    $ws.Cells.Item($row, 6).Value = $planDur
    $ws.Cells.Item($row, 7).Value = $actualStart
    $ws.Cells.Item($row, 8).Value = $actualDur
    $ws.Cells.Item($row, 9).Value = $percent
}

# 2.2 Manually add a new student to a class.
Set-PlanActual $projFunctions 16 25 10 25 10 1
# 2.3 Edit an existing student.
Set-PlanActual $projFunctions 17 25 10 25 10 1
# 2.4 Remove a student.
Set-PlanActual $projFunctions 18 30 10 30 10 1
# 2.5 Change students from class A to class B
Set-PlanActual $projFunctions 19 32 10 32 0 0
# 2.6 View list of classes.
Set-PlanActual $projFunctions 20 32 10 32 10 1
# 2.7 View list of students in a class.
Set-PlanActual $projFunctions 21 32 10 32 10 1

# 3.4 Edit an existing course. (editCourse)
Set-PlanActual $projFunctions 26 32 7 32 7 1
# 3.5 Remove a course. (removeCourse)
Set-PlanActual $projFunctions 27 32 7 32 7 1

# Row 32 is function 3.10 ("View list of courses in the current
# semester."); format it with 2 decimals so it reads "3.10".
$projFunctions.Cells.Item(32, 1).NumberFormat = "0.00"

# 3.11 Create / update / delete / view all lecturers. -> assigned to Trí
$projFunctions.Cells.Item(33, 4).Value = "Trí"
Set-PlanActual $projFunctions 33 41 2 41 2 1

# 7.2 View check-in result.
Set-PlanActual $projFunctions 50 25 7 27 2 1
# 7.3 View schedules.
Set-PlanActual $projFunctions 51 41 2 41 2 1
# 7.4 View his/her scores of a course.
Set-PlanActual $projFunctions 52 25 7 27 2 1

# Widen column A (code numbers) from 6 to 7 characters.
$projFunctions.Columns.Item(1).ColumnWidth = 6.285714285714286

# ---------------------------------------------------------------------
# Sheet "NoteDetails": add a "Group's Notes" section documenting reused
# code between similar functions.
# ---------------------------------------------------------------------
$noteDetails = $wb.Worksheets.Item("NoteDetails")
$noteDetails.Activate()

$noteDetails.Cells.Item(12, 1).Value = "Group's Notes"

# Populate in this order so new shared strings land at the same table
# offsets as the authored workbook.
$noteDetails.Cells.Item(14, 3).Value = "Same as 6,2."
$noteDetails.Cells.Item(13, 3).Value = "Same as 6,1."
$noteDetails.Cells.Item(15, 3).Value = "Same as 6,3. Reuse codes from 5,1."
$noteDetails.Cells.Item(16, 3).Value = "Reuse codes from 4,1."

$noteDetails.Cells.Item(13, 1).Value = 3.8
$noteDetails.Cells.Item(14, 1).Value = 3.9
$noteDetails.Cells.Item(15, 1).Value = 3.1
$noteDetails.Cells.Item(15, 1).NumberFormat = "0.00"
$noteDetails.Cells.Item(16, 1).Value = 6.7

$noteDetails.Range("C17").Select()

# Restore ProjFunctions as the active sheet/selection.
$projFunctions.Activate()
$projFunctions.Range("I20").Select()
